# Generate Report for Handoff
# The localization handoff transform failed, so:
#  - Status changes from "Ready for handoff" to "Handoff transform failed"
#  - The Latest Handoff File (and its hyperlink) is cleared
#  - The Latest Handoff Datetime / Latest Handback DateTime reset to the zero-date
#  - The Handoff Reason becomes "Ignored" for both rows

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the Status text everywhere it currently reads "Ready for handoff"
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"
$zhcn.Range("B2").Value = "Handoff transform failed"
$dede.Range("B2").Value = "Handoff transform failed"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/5e11de6b9e7f3f72068f0d59985cb9e7bde843b0/e2e/88e15722-ee89-4119-8452-37afc5243a5f.md"
$cfgAddress = "https://github.com/OpenLocalizationTest/oltest/blob/5e11de6b9e7f3f72068f0d59985cb9e7bde843b0/.localization-config"

foreach ($ws in @($zhcn, $dede)) {
    # This runtime only supports removing a single cell's hyperlink by
    # clearing the whole sheet hyperlink collection and re-adding the ones
    # that should survive (A2 and A3); the C2 "Latest Handoff File" link is
    # intentionally not recreated below.
    $ws.Hyperlinks.Delete()

    # Remove the "Latest Handoff File" value (the handoff transform failed,
    # so there is no handoff file to report any more)
    $ws.Range("C2").Clear()

    # Reset datetimes to the epoch placeholder value
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"

    # Handoff Reason both become "Ignored"
    $ws.Range("H2").Value = "Ignored"
    $ws.Range("H3").Value = "Ignored"

    # Restore the two hyperlinks that are still valid
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, "88e15722-ee89-4119-8452-37afc5243a5f.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $cfgAddress, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
}
